# This script reapplies the per-Leve profit recompute captured in the
# upstream diff for Sheets/Marilith_Profits.xlsx. Each FFXIV crafting job
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) lives on its own worksheet; the columns
# H:N hold the scheduled-runner's market-price/profit recompute
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]).
# Blank cells in the 'Value' map mean the cell is cleared (no longer present).

$wb = $excel.ActiveWorkbook

$edits = @(
    @{ Sheet="ALC"; Row=2; Cells=@{ 8=1282.2941; 9=116.583336; 10=4080; 11=116.583336; 12=4080; 13=-3.583336000000003; 14=-4306 } }
    @{ Sheet="ALC"; Row=12; Cells=@{ 8=383.33334; 9=320; 10=700; 11=320; 12=700; 13=-150; 14=-1040 } }
    @{ Sheet="ALC"; Row=28; Cells=@{ 8=6987.1333; 9=7475.7856; 11=7475.7856; 13=-6990.7856 } }
    @{ Sheet="ALC"; Row=29; Cells=@{ 8=135.25; 9=145.5; 10=125; 11=436.5; 12=375; 13=-155.5; 14=-937 } }
    @{ Sheet="ALC"; Row=43; Cells=@{ 8=750000; 9=0; 10=750000; 11=0; 12=750000; 13=""; 14=-750138 } }
    @{ Sheet="ALC"; Row=58; Cells=@{ 8=1221.7273; 10=1623.1666; 12=4869.4998; 14=-5169.4998 } }
    @{ Sheet="ALC"; Row=64; Cells=@{ 8=7000; 9=7000; 11=7000; 13=-6752 } }
    @{ Sheet="ALC"; Row=67; Cells=@{ 8=7000; 9=7000; 11=7000; 13=-6142 } }
    @{ Sheet="ALC"; Row=92; Cells=@{ 8=258.08334; 9=266.77777; 11=266.77777; 13=981.2222300000001 } }
    @{ Sheet="ALC"; Row=112; Cells=@{ 8=1849.2941; 9=1365.5; 10=2113.182; 11=4096.5; 12=6339.545999999999; 13=-2988.5; 14=-8555.545999999998 } }
    @{ Sheet="ALC"; Row=125; Cells=@{ 8=6979.8; 9=3633.3333; 11=32699.9997; 13=-30239.9997 } }
    @{ Sheet="ALC"; Row=135; Cells=@{ 8=1136.3077; 9=981; 11=8829; 13=-6294 } }
    @{ Sheet="ALC"; Row=137; Cells=@{ 8=2375.7144; 9=2124.5; 11=6373.5; 13=-3823.5 } }
    @{ Sheet="ALC"; Row=138; Cells=@{ 8=1668.4375; 9=939; 10=2000; 11=2817; 12=6000; 13=2323; 14=-16280 } }
    @{ Sheet="ARM"; Row=5; Cells=@{ 8=28.625; 9=30; 11=30; 13=82 } }
    @{ Sheet="ARM"; Row=32; Cells=@{ 8=6247.654; 9=5483.2915; 11=5483.2915; 13=-5196.2915 } }
    @{ Sheet="ARM"; Row=55; Cells=@{ 8=28333.334; 9=0; 10=28333.334; 11=0; 12=28333.334; 13=""; 14=-28963.334 } }
    @{ Sheet="ARM"; Row=109; Cells=@{ 8=46666.668; 10=46666.668; 12=46666.668; 14=-49440.668 } }
    @{ Sheet="BSM"; Row=4; Cells=@{ 8=28.625; 9=30; 11=30; 13=85 } }
    @{ Sheet="BSM"; Row=22; Cells=@{ 8=860.2; 9=860.2; 11=860.2; 13=-687.2 } }
    @{ Sheet="BSM"; Row=135; Cells=@{ 8=43499.75; 10=43499.75; 12=43499.75; 14=-53639.75 } }
    @{ Sheet="BSM"; Row=137; Cells=@{ 8=74332; 10=74332; 12=74332; 14=-84532 } }
    @{ Sheet="CRP"; Row=62; Cells=@{ 8=3150; 10=2975; 12=2975; 14=-4223 } }
    @{ Sheet="CRP"; Row=65; Cells=@{ 8=3150; 10=2975; 12=14875; 14=-21115 } }
    @{ Sheet="CRP"; Row=99; Cells=@{ 8=4291.857; 9=4790.5; 10=1300; 11=4790.5; 12=1300; 13=-3292.5; 14=-4296 } }
    @{ Sheet="CRP"; Row=103; Cells=@{ 8=3549.25; 9=3549.25; 11=3549.25; 13=-2377.25 } }
    @{ Sheet="CRP"; Row=126; Cells=@{ 8=4291.857; 9=4790.5; 10=1300; 11=14371.5; 12=3900; 13=-11901.5; 14=-8840 } }
    @{ Sheet="GSM"; Row=43; Cells=@{ 8=25539.572; 10=32155.4; 12=32155.4; 14=-32457.4 } }
    @{ Sheet="GSM"; Row=122; Cells=@{ 8=8336166; 9=9618038; 11=28854114; 13=-28851664 } }
    @{ Sheet="LTW"; Row=16; Cells=@{ 8=3447.7144; 9=1949.25; 10=5445.6665; 11=1949.25; 12=5445.6665; 13=-1779.25; 14=-5785.6665 } }
    @{ Sheet="LTW"; Row=22; Cells=@{ 8=1298.5333; 9=1085.7142; 11=1085.7142; 13=-790.7141999999999 } }
    @{ Sheet="LTW"; Row=27; Cells=@{ 8=1298.5333; 9=1085.7142; 11=1085.7142; 13=-978.7141999999999 } }
    @{ Sheet="LTW"; Row=46; Cells=@{ 8=2654.9285; 9=2713.3333; 10=2549.8; 11=2713.3333; 12=2549.8; 13=-2525.3333; 14=-2925.8 } }
    @{ Sheet="LTW"; Row=100; Cells=@{ 8=1337.6; 9=922; 11=922; 13=-381 } }
    @{ Sheet="LTW"; Row=109; Cells=@{ 8=69990; 10=69990; 12=69990; 14=-72764 } }
    @{ Sheet="LTW"; Row=122; Cells=@{ 8=3473.6365; 10=4333; 12=12999; 14=-17899 } }
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    foreach ($col in $edit.Cells.Keys) {
        $ws.Cells.Item($edit.Row, $col).Value = $edit.Cells[$col]
    }
}